$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current (last data) row 76, pushing the
# existing row 76 down to row 77. The new row 76 carries the newer
# weekly price observation for Terminal Hortofrutícola Agro Chillán - Mango.
$ws.Rows.Item(76).Insert()

$ws.Range("A76").Value = 7
$ws.Range("B76").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C76").Value = "Ñuble"
$ws.Range("D76").Value = 44656
$ws.Range("E76").Value = 16
$ws.Range("F76").Value = "Fruta"
$ws.Range("G76").Value = 100108
$ws.Range("H76").Value = "Tropicales y subtropicales"
$ws.Range("I76").Value = 100108002
$ws.Range("J76").Value = "Mango"
$ws.Range("K76").Value = "Sin especificar"
$ws.Range("L76").Value = "Primera"
$ws.Range("M76").Value = 60
$ws.Range("N76").Value = 7500
$ws.Range("O76").Value = 8000
$ws.Range("P76").Value = 7750
$ws.Range("Q76").Value = "$/bandeja 4 kilos"
$ws.Range("R76").Value = "Perú"
$ws.Range("S76").Value = 1938
$ws.Range("T76").Value = 4
